$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value that was bumped by one day
# (2023-10-06 -> 2023-10-07, serial 45205 -> 45206) for every data row.
# Data rows run from row 2 through row 498.
$ws.Range("C2:C498").Value = 45206
